$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($ws, $ref, $val)
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-CellText $ws "D2" "43.108.13"
Set-CellText $ws "E2" "  +0.21%  "
Set-CellText $ws "D3" "2.314.16"
Set-CellText $ws "E3" "  +0.27%  "
Set-CellText $ws "E4" "  +0.02%  "
Set-CellText $ws "D5" "302.36"
Set-CellText $ws "E5" "  -0.16%  "
Set-CellText $ws "D6" "98.97"
Set-CellText $ws "E6" "  -2.19%  "
Set-CellText $ws "D7" "0.521"
Set-CellText $ws "E7" "  +3.03%  "
Set-CellText $ws "E8" "  +0.02%  "
Set-CellText $ws "D9" "0.522"
Set-CellText $ws "E9" "  +0.65%  "
Set-CellText $ws "D10" "35.83"
Set-CellText $ws "E10" "  +1.28%  "
Set-CellText $ws "D11" "0.0790"
Set-CellText $ws "E11" "  -0.66%  "
Set-CellText $ws "E12" "  -0.91%  "
Set-CellText $ws "E13" "  -0.74%  "
Set-CellText $ws "D14" "6.94"
Set-CellText $ws "E14" "  +0.36%  "
Set-CellText $ws "D15" "2.673.73"
Set-CellText $ws "E15" "  -0.38%  "
Set-CellText $ws "D16" "2.264.69"
Set-CellText $ws "E16" "  -1.06%  "
Set-CellText $ws "D17" "0.791"
Set-CellText $ws "D18" "43.018.82"
Set-CellText $ws "E18" "  +0.21%  "
Set-CellText $ws "D19" "13.54"
Set-CellText $ws "E19" "  +7.28%  "
Set-CellText $ws "D20" "0.0₃0911"
Set-CellText $ws "E20" "  +0.61%  "
Set-CellText $ws "E21" "  +0.04%  "
Set-CellText $ws "D22" "68.15"
Set-CellText $ws "E22" "  +0.29%  "
Set-CellText $ws "D23" "240.37"
Set-CellText $ws "E23" "  +1.30%  "
Set-CellText $ws "E24" "  -1.25%  "
Set-CellText $ws "B25" "PancakeSwap"
Set-CellText $ws "C25" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-CellText $ws "D25" "2.46"
Set-CellText $ws "E25" "  -0.35%  "
Set-CellText $ws "B26" "Dai"
Set-CellText $ws "C26" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-CellText $ws "D26" "0.999"
Set-CellText $ws "E26" "  -0.05%  "
Set-CellText $ws "D27" "24.96"
Set-CellText $ws "E27" "  +0.63%  "
Set-CellText $ws "D28" "168.55"
Set-CellText $ws "E28" "  +0.38%  "
Set-CellText $ws "B29" "Toncoin"
Set-CellText $ws "C29" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-CellText $ws "D29" "2.05"
Set-CellText $ws "E29" "  -6.16%  "
Set-CellText $ws "B30" "Cosmos"
Set-CellText $ws "C30" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-CellText $ws "D30" "9.19"
Set-CellText $ws "E30" "  -0.67%  "
Set-CellText $ws "D31" "33.41"
Set-CellText $ws "E31" "  -2.20%  "
Set-CellText $ws "B32" "Filecoin"
Set-CellText $ws "C32" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-CellText $ws "D32" "5.23"
Set-CellText $ws "E32" "  +4.01%  "
Set-CellText $ws "B33" "RenderToken"
Set-CellText $ws "C33" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-CellText $ws "D33" "4.92"
Set-CellText $ws "E33" "  +6.37%  "
Set-CellText $ws "D34" "18.48"
Set-CellText $ws "E34" "  +8.71%  "
Set-CellText $ws "E36" "  -0.08%  "
Set-CellText $ws "D37" "0.0695"
Set-CellText $ws "E37" "  +0.47%  "
Set-CellText $ws "B38" "ARBITRUM"
Set-CellText $ws "C38" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-CellText $ws "D38" "1.81"
Set-CellText $ws "E38" "  +1.06%  "
Set-CellText $ws "B39" "Kaspa"
Set-CellText $ws "C39" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-CellText $ws "D39" "0.102"
Set-CellText $ws "E39" "  +0.41%  "
Set-CellText $ws "B40" "Stellar"
Set-CellText $ws "C40" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-CellText $ws "D40" "0.111"
Set-CellText $ws "E40" "  +0.96%  "
Set-CellText $ws "B41" "LidoDAOToken"
Set-CellText $ws "C41" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-CellText $ws "D41" "2.77"
Set-CellText $ws "E41" "  -2.09%  "
Set-CellText $ws "D42" "1.998.07"
Set-CellText $ws "E42" "  -0.23%  "
Set-CellText $ws "E43" "  +0.17%  "
Set-CellText $ws "E44" "  -6.22%  "
Set-CellText $ws "D45" "10.11"
Set-CellText $ws "E45" "  -1.57%  "
Set-CellText $ws "D46" "17.49"
Set-CellText $ws "E46" "  -0.53%  "
Set-CellText $ws "E47" "  -0.85%  "
Set-CellText $ws "D48" "54.89"
Set-CellText $ws "E48" "  -1.61%  "
Set-CellText $ws "D49" "74.33"
Set-CellText $ws "E49" "  +5.82%  "
Set-CellText $ws "D50" "2.539.73"
Set-CellText $ws "E50" "  +0.85%  "
Set-CellText $ws "E51" "  +1.20%  "
